$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "ValidLogin"

# Add Password / manager columns
$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "manager"

# Update selection to B3
$ws.Range("B3").Select()
